$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 37.166668
$ws.Range("I6").Value = 37.8
$ws.Range("J6").Value = 34
$ws.Range("K6").Value = 113.4
$ws.Range("L6").Value = 102
$ws.Range("M6").Value = -1.399999999999991
$ws.Range("N6").Value = -326
$ws.Range("H11").Value = 29067.5
$ws.Range("I11").Value = 29067.5
$ws.Range("K11").Value = 29067.5
$ws.Range("M11").Value = -28927.5
$ws.Range("H12").Value = 132.42857
$ws.Range("I12").Value = 133
$ws.Range("K12").Value = 133
$ws.Range("M12").Value = 37
$ws.Range("H19").Value = 509.30435
$ws.Range("I19").Value = 631.0909
$ws.Range("J19").Value = 397.66666
$ws.Range("K19").Value = 631.0909
$ws.Range("L19").Value = 397.66666
$ws.Range("M19").Value = -456.0909
$ws.Range("N19").Value = -747.66666
$ws.Range("H51").Value = 3315.5454
$ws.Range("I51").Value = 2970.5881
$ws.Range("J51").Value = 4488.4
$ws.Range("K51").Value = 2970.5881
$ws.Range("L51").Value = 4488.4
$ws.Range("M51").Value = -2486.5881
$ws.Range("N51").Value = -5456.4
$ws.Range("H104").Value = 631.3333
$ws.Range("I104").Value = 631.3333
$ws.Range("K104").Value = 1893.9999
$ws.Range("M104").Value = -146.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 627.6875
$ws.Range("I2").Value = 564.9231
$ws.Range("K2").Value = 564.9231
$ws.Range("M2").Value = -451.9231
$ws.Range("H97").Value = 1733.5454
$ws.Range("I97").Value = 906.95
$ws.Range("K97").Value = 906.95
$ws.Range("M97").Value = -410.95
$ws.Range("H116").Value = 627.6875
$ws.Range("I116").Value = 564.9231
$ws.Range("K116").Value = 564.9231
$ws.Range("M116").Value = 1729.0769
$ws.Range("H132").Value = 1391244
$ws.Range("I132").Value = 2085257.6
$ws.Range("J132").Value = 3216.5
$ws.Range("K132").Value = 6255772.800000001
$ws.Range("L132").Value = 9649.5
$ws.Range("M132").Value = -6253242.800000001
$ws.Range("N132").Value = -14709.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 627.6875
$ws.Range("I3").Value = 564.9231
$ws.Range("K3").Value = 564.9231
$ws.Range("M3").Value = -450.9231
$ws.Range("H86").Value = 7000
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 9000
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -11246
$ws.Range("H89").Value = 7000
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 9000
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -56232
$ws.Range("H94").Value = 4275
$ws.Range("I94").Value = 3756.8667
$ws.Range("K94").Value = 3756.8667
$ws.Range("M94").Value = -3305.8667
$ws.Range("H105").Value = 2500.2888
$ws.Range("I105").Value = 2198.2903
$ws.Range("K105").Value = 2198.2903
$ws.Range("M105").Value = -451.2903000000001
$ws.Range("H107").Value = 958.0909
$ws.Range("I107").Value = 868.8333
$ws.Range("K107").Value = 868.8333
$ws.Range("M107").Value = 1051.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 171.85715
$ws.Range("I7").Value = 132.54546
$ws.Range("J7").Value = 215.1
$ws.Range("K7").Value = 132.54546
$ws.Range("L7").Value = 215.1
$ws.Range("M7").Value = -19.54545999999999
$ws.Range("N7").Value = -441.1
$ws.Range("H11").Value = 505
$ws.Range("I11").Value = 837.6667
$ws.Range("K11").Value = 837.6667
$ws.Range("M11").Value = -697.6667
$ws.Range("H48").Value = 10000
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 10000
$ws.Range("M48").Value = -9524

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 58.153847
$ws.Range("J2").Value = 55.8
$ws.Range("L2").Value = 334.8
$ws.Range("N2").Value = -560.8
$ws.Range("H7").Value = 58
$ws.Range("J7").Value = 63.4
$ws.Range("L7").Value = 190.2
$ws.Range("N7").Value = -414.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4574.5713
$ws.Range("I18").Value = 2003.6666
$ws.Range("K18").Value = 2003.6666
$ws.Range("M18").Value = -1710.6666
$ws.Range("H36").Value = 2000
$ws.Range("I36").Value = 2000
$ws.Range("K36").Value = 2000
$ws.Range("M36").Value = -1515
$ws.Range("H43").Value = 9198.4
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 14060
$ws.Range("J46").Value = 18733
$ws.Range("L46").Value = 18733
$ws.Range("N46").Value = -19045
$ws.Range("H55").Value = 9999.666999999999
$ws.Range("J55").Value = 14499.5
$ws.Range("L55").Value = 14499.5
$ws.Range("N55").Value = -15153.5
$ws.Range("H57").Value = 24499.666
$ws.Range("J57").Value = 24499.666
$ws.Range("L57").Value = 24499.666
$ws.Range("N57").Value = -26139.666
$ws.Range("H70").Value = 39719.785
$ws.Range("I70").Value = 19411.143
$ws.Range("K70").Value = 19411.143
$ws.Range("M70").Value = -19141.143
$ws.Range("H73").Value = 39719.785
$ws.Range("I73").Value = 19411.143
$ws.Range("K73").Value = 19411.143
$ws.Range("M73").Value = -18475.143
$ws.Range("H126").Value = 4738.9644
$ws.Range("I126").Value = 5355.857
$ws.Range("J126").Value = 2888.2856
$ws.Range("K126").Value = 16067.571
$ws.Range("L126").Value = 8664.856800000001
$ws.Range("M126").Value = -13597.571
$ws.Range("N126").Value = -13604.8568
$ws.Range("H132").Value = 13776.881
$ws.Range("I132").Value = 16404.344
$ws.Range("K132").Value = 49213.03200000001
$ws.Range("M132").Value = -46683.03200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2651.04
$ws.Range("J22").Value = 3309.889
$ws.Range("L22").Value = 3309.889
$ws.Range("N22").Value = -3899.889
$ws.Range("H27").Value = 2651.04
$ws.Range("J27").Value = 3309.889
$ws.Range("L27").Value = 3309.889
$ws.Range("N27").Value = -3523.889
$ws.Range("H46").Value = 5005.8
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 5291.9287
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 5291.9287
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -5667.9287
$ws.Range("H55").Value = 1231
$ws.Range("I55").Value = 1386
$ws.Range("K55").Value = 1386
$ws.Range("M55").Value = -1213
$ws.Range("H81").Value = 49900
$ws.Range("J81").Value = 49900
$ws.Range("L81").Value = 49900
$ws.Range("N81").Value = -51896
$ws.Range("H84").Value = 49900
$ws.Range("J84").Value = 49900
$ws.Range("L84").Value = 149700
$ws.Range("N84").Value = -159684

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 135.90909
$ws.Range("I2").Value = 135.90909
$ws.Range("K2").Value = 135.90909
$ws.Range("M2").Value = -23.90908999999999
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 157.14285
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 157.14285
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -44.14285000000001
$ws.Range("N4").Value = -326
$ws.Range("H81").Value = 3063.9443
$ws.Range("I81").Value = 2080.182
$ws.Range("J81").Value = 4609.857
$ws.Range("K81").Value = 4160.364
$ws.Range("L81").Value = 9219.714
$ws.Range("M81").Value = -3099.364
$ws.Range("N81").Value = -11341.714
$ws.Range("H84").Value = 3063.9443
$ws.Range("I84").Value = 2080.182
$ws.Range("J84").Value = 4609.857
$ws.Range("K84").Value = 20801.82
$ws.Range("L84").Value = 46098.57
$ws.Range("M84").Value = -15497.82
$ws.Range("N84").Value = -56706.57
